# Update LTP/PREV figures on the "ltp" worksheet with the latest modular
# amount-based calc values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ltp")

$ws.Range("B2").Value = 2233.35
$ws.Range("C2").Value = 2246
$ws.Range("B3").Value = 384.1
$ws.Range("C3").Value = 393
$ws.Range("B4").Value = 1569.8
$ws.Range("C4").Value = 1564.4
$ws.Range("B5").Value = 7495.6
$ws.Range("C5").Value = 7562.05
$ws.Range("B6").Value = 216.1
$ws.Range("C6").Value = 217.85
$ws.Range("B7").Value = 191.25
$ws.Range("C7").Value = 195.6
$ws.Range("B8").Value = 43926.85
$ws.Range("C8").Value = 43839.8
$ws.Range("B9").Value = 592.45
$ws.Range("C9").Value = 592
$ws.Range("B10").Value = 3509.05
$ws.Range("C10").Value = 3510.95
$ws.Range("B11").Value = 147.3
$ws.Range("C11").Value = 144.1
$ws.Range("B12").Value = 1273.95
$ws.Range("C12").Value = 1269.7
$ws.Range("B13").Value = 1486.4
$ws.Range("C13").Value = 1470.25
$ws.Range("B14").Value = 630.35
$ws.Range("C14").Value = 622.2
$ws.Range("B15").Value = 461
$ws.Range("C15").Value = 464.85
$ws.Range("B16").Value = 1479.75
$ws.Range("C16").Value = 1488.3
$ws.Range("B17").Value = 263.95
$ws.Range("C17").Value = 262.7
$ws.Range("B18").Value = 19476.9
$ws.Range("C18").Value = 19486.5
$ws.Range("B19").Value = 579.75
$ws.Range("C19").Value = 574.35
$ws.Range("B20").Value = 652.3
$ws.Range("C20").Value = 651.45
$ws.Range("B21").Value = 645
$ws.Range("C21").Value = 646.05
$ws.Range("B22").Value = 249.8
$ws.Range("C22").Value = 249.55
$ws.Range("B23").Value = 119.65
$ws.Range("C23").Value = 119.4
